{"js": "// Update the \"Sprint No.\" value from 1 -> 2 and the \"Review Date\" value\n// from 02/09/18 -> 02/21/18 in the header table of the code-review\n// checklist document. Both values live in the first table on the page.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst headerTable = tables.items[0];\nconst rows = headerTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 1 (\"Reviewer's Name\" | name | \"Sprint No.\" | value) -> 4th cell holds \"1\"\nconst sprintRow = rows.items[1];\nconst sprintCells = sprintRow.cells;\nsprintCells.load(\"items\");\nawait context.sync();\n\nconst sprintValueCell = sprintCells.items[3];\nconst sprintMatches = sprintValueCell.body.search(\"1\", { matchCase: true, matchWholeWord: true });\nsprintMatches.load(\"items\");\nawait context.sync();\nsprintMatches.items[0].insertText(\"2\", Word.InsertLocation.replace);\n\n// Row 2 (\"Review Date\" | value) -> 2nd cell holds \"02/09/18\"\nconst dateRow = rows.items[2];\nconst dateCells = dateRow.cells;\ndateCells.load(\"items\");\nawait context.sync();\n\nconst dateValueCell = dateCells.items[1];\nconst dateMatches = dateValueCell.body.search(\"02/09/18\", { matchCase: true });\ndateMatches.load(\"items\");\nawait context.sync();\ndateMatches.items[0].insertText(\"02/21/18\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Update the \"Sprint No.\" value from 1 -> 2 and the \"Review Date\" value\n# from 02/09/18 -> 02/21/18 in the header table of the code-review\n# checklist document (first table on the page).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row 2, Column 4 holds the Sprint No. value (\"1\"). Scope the Find to the\n# cell's own Range so only this occurrence of \"1\" is touched, then set the\n# (now narrowed-to-the-match) range's Text directly -- this preserves the\n# surrounding run/paragraph formatting instead of replacing the whole run.\n$sprintCell = $t.Cell(2, 4)\n$sprintRng = $sprintCell.Range\n$sprintFind = $sprintRng.Find\n$sprintFind.ClearFormatting()\n$sprintFound = $sprintFind.Execute(\"1\")\nif ($sprintFound) { $sprintRng.Text = \"2\" }\n\n# Row 3, Column 2 holds the Review Date value (\"02/09/18\").\n$dateCell = $t.Cell(3, 2)\n$dateRng = $dateCell.Range\n$dateFind = $dateRng.Find\n$dateFind.ClearFormatting()\n$dateFound = $dateFind.Execute(\"02/09/18\")\nif ($dateFound) { $dateRng.Text = \"02/21/18\" }\n"}
